$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value  = "QA correcto"
$ws.Range("B3").Value  = "QA incorrecto"
$ws.Range("B4").Value  = "QA Reporte"
$ws.Range("B5").Value  = "Descargar reporte"
$ws.Range("B6").Value  = "El reporte de QA no ha encontrado errores en los datos."
$ws.Range("B7").Value  = "El reporte de QA ha encontrado errores en los datos, revisar detalles en el informe de reporte."
$ws.Range("B8").Value  = "Configuración correcta"
$ws.Range("B9").Value  = "Configuración incorrecta"
$ws.Range("B10").Value = "Todos los paquetes fueron instalados"
$ws.Range("B11").Value = "Algunos paquetes no fueron instalados:"
$ws.Range("B12").Value = "Cerrar"
